$p = $ppt.ActivePresentation
$s5 = $p.Slides.Item(5)
$sh = $s5.Shapes.Item(2)
$tbl = $sh.Table

$guids = @(
  "{5940675A-B579-460E-94D1-54222C63F5DA}",
  "{3C2FFA5D-87B4-456A-9821-1D502468CF0F}",
  "{3AE1B340-825A-4F07-AC60-6A056DE400A1}",
  "{2D5ABB26-0587-4C30-8999-92F81FD0307C}"
)
foreach ($g in $guids) {
    $tbl.Style = $g
    Write-Host "tried $g -> now:" $tbl.Style
}
